$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 20:15"

# Update country statistics that changed between the two snapshots.
# (Country name labels in column A are unchanged; only the numeric
# columns B (Casos totales), C (Nuevos casos), D (Casos activos),
# E (Recuperados), G (Muertes hoy) and H (Muertes) are refreshed.)
# Row 4
$ws.Range("B4").Value = 6359512
$ws.Range("C4").Value = 24268
$ws.Range("D4").Value = 3584207
$ws.Range("E4").Value = 2583827
$ws.Range("G4").Value = 420
$ws.Range("H4").Value = 191478

# Row 19
$ws.Range("B19").Value = 309156
$ws.Range("C19").Value = 8975
$ws.Range("E19").Value = 191244

# Row 24
$ws.Range("B24").Value = 249569
$ws.Range("C24").Value = 755
$ws.Range("E24").Value = 15569
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 9400

# Row 32
$ws.Range("B32").Value = 117175
$ws.Range("C32").Value = 815
$ws.Range("D32").Value = 102200
$ws.Range("E32").Value = 8301
$ws.Range("G32").Value = 26
$ws.Range("H32").Value = 6674

# Row 49
$ws.Range("B49").Value = 68605
$ws.Range("C49").Value = 1750
$ws.Range("D49").Value = 52483
$ws.Range("E49").Value = 14830
$ws.Range("G49").Value = 39
$ws.Range("H49").Value = 1292

# Row 53
$ws.Range("B53").Value = 56516
$ws.Range("C53").Value = 1303
$ws.Range("D53").Value = 20612
$ws.Range("E53").Value = 35024
$ws.Range("G53").Value = 24
$ws.Range("H53").Value = 880

# Row 93
$ws.Range("B93").Value = 10149
$ws.Range("C93").Value = 24
$ws.Range("E93").Value = 523

# Row 100
$ws.Range("B100").Value = 8361
$ws.Range("C100").Value = 80
$ws.Range("D100").Value = 5642
$ws.Range("E100").Value = 2690
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 29

# Row 101
$ws.Range("B101").Value = 8323
$ws.Range("C101").Value = 241
$ws.Range("D101").Value = 3611
$ws.Range("E101").Value = 4625
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 87

# Row 102
$ws.Range("B102").Value = 8301
$ws.Range("D102").Value = 5870
$ws.Range("E102").Value = 2221
$ws.Range("H102").Value = 210

# Row 113
$ws.Range("B113").Value = 4780
$ws.Range("C113").Value = 60
$ws.Range("D113").Value = 3844
$ws.Range("E113").Value = 842
$ws.Range("H113").Value = 94

# Row 114
$ws.Range("B114").Value = 4729
$ws.Range("D114").Value = 1807
$ws.Range("E114").Value = 2860
$ws.Range("H114").Value = 62

# Row 119
$ws.Range("B119").Value = 4265
$ws.Range("C119").Value = 58
$ws.Range("D119").Value = 2511
$ws.Range("E119").Value = 1728
$ws.Range("H119").Value = 26

# Row 120
$ws.Range("B120").Value = 4255
$ws.Range("D120").Value = 2163
$ws.Range("E120").Value = 2074
$ws.Range("H120").Value = 18

# Row 121
$ws.Range("B121").Value = 4215
$ws.Range("D121").Value = 3318
$ws.Range("E121").Value = 824
$ws.Range("H121").Value = 73

# Row 122
$ws.Range("B122").Value = 4214
$ws.Range("D122").Value = 2370
$ws.Range("E122").Value = 640
$ws.Range("H122").Value = 100

# Row 136
$ws.Range("B136").Value = 2536
$ws.Range("C136").Value = 3
$ws.Range("E136").Value = 1199

# Row 146
$ws.Range("C146").Value = 278
$ws.Range("D146").Value = 493
$ws.Range("E146").Value = 1501
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 8

# Row 147
$ws.Range("B147").Value = 2002
$ws.Range("C147").Value = 90
$ws.Range("D147").Value = 880
$ws.Range("E147").Value = 1112
$ws.Range("H147").Value = 10

# Row 148
$ws.Range("B148").Value = 1984
$ws.Range("C148").Value = 19
$ws.Range("D148").Value = 1565
$ws.Range("E148").Value = 406
$ws.Range("H148").Value = 13

# Row 149
$ws.Range("B149").Value = 1983
$ws.Range("C149").Value = 4
$ws.Range("D149").Value = 1194
$ws.Range("E149").Value = 217
$ws.Range("H149").Value = 572

# Row 150
$ws.Range("B150").Value = 1764
$ws.Range("C150").Value = 5
$ws.Range("D150").Value = 1630
$ws.Range("E150").Value = 111
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 23

# Row 153
$ws.Range("B153").Value = 1502
$ws.Range("C153").Value = 4
$ws.Range("E153").Value = 342

# Row 162
$ws.Range("B162").Value = 1148
$ws.Range("C162").Value = 63
$ws.Range("D162").Value = 528
$ws.Range("E162").Value = 589
$ws.Range("H162").Value = 31

# Row 163
$ws.Range("B163").Value = 1118
$ws.Range("D163").Value = 255
$ws.Range("E163").Value = 850
$ws.Range("H163").Value = 13

# Row 164
$ws.Range("B164").Value = 1111
$ws.Range("D164").Value = 359
$ws.Range("E164").Value = 746
$ws.Range("H164").Value = 6
